$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Nitish Rana" in the source data is followed by a non-breaking space
# (U+00A0, matching the existing rows) rather than a regular space.
$nbsp = [string]([char]0x00A0)

# New match-log rows scraped for Nitish Rana (IPL 2020 UAE season),
# appended below the existing table (previously A1:K15, now A1:K29).
$newRows = @(
    ,@(" Sharjah", " October 26 2020", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kolkata Knight Riders", "Kings XI Punjab", ("Nitish Rana" + $nbsp), "0", "1", "0", "0", "0.00")
    ,@(" Abu Dhabi", " September 23 2020", "Mumbai won by 49 runs", "Kolkata Knight Riders", "Mumbai Indians", ("Nitish Rana" + $nbsp), "24", "18", "2", "1", "133.33")
    ,@(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Kolkata Knight Riders", "Sunrisers Hyderabad", ("Nitish Rana" + $nbsp), "29", "20", "3", "1", "145.00")
    ,@(" Dubai (DSC)", " October 29 2020", "Super Kings won by 6 wickets", "Kolkata Knight Riders", "Chennai Super Kings", ("Nitish Rana" + $nbsp), "87", "61", "10", "4", "142.62")
    ,@(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Kolkata Knight Riders", "Delhi Capitals", ("Nitish Rana" + $nbsp), "81", "53", "13", "1", "152.83")
    ,@(" Abu Dhabi", " September 26 2020", "KKR won by 7 wickets (with 12 balls remaining)", "Kolkata Knight Riders", "Sunrisers Hyderabad", ("Nitish Rana" + $nbsp), "26", "13", "6", "0", "200.00")
    ,@(" Dubai (DSC)", " September 30 2020", "KKR won by 37 runs", "Kolkata Knight Riders", "Rajasthan Royals", ("Nitish Rana" + $nbsp), "22", "17", "2", "1", "129.41")
    ,@(" Abu Dhabi", " October 16 2020", "Mumbai won by 8 wickets (with 19 balls remaining)", "Kolkata Knight Riders", "Mumbai Indians", ("Nitish Rana" + $nbsp), "5", "6", "1", "0", "83.33")
    ,@(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", ("Nitish Rana" + $nbsp), "0", "1", "0", "0", "0.00")
    ,@(" Abu Dhabi", " October 10 2020", "KKR won by 2 runs", "Kolkata Knight Riders", "Kings XI Punjab", ("Nitish Rana" + $nbsp), "2", "4", "0", "0", "50.00")
    ,@(" Sharjah", " October 03 2020", "Capitals won by 18 runs", "Kolkata Knight Riders", "Delhi Capitals", ("Nitish Rana" + $nbsp), "58", "35", "4", "4", "165.71")
    ,@(" Abu Dhabi", " October 21 2020", "RCB won by 8 wickets (with 39 balls remaining)", "Kolkata Knight Riders", "Royal Challengers Bangalore", ("Nitish Rana" + $nbsp), "0", "1", "0", "0", "0.00")
    ,@(" Sharjah", " October 12 2020", "RCB won by 82 runs", "Kolkata Knight Riders", "Royal Challengers Bangalore", ("Nitish Rana" + $nbsp), "9", "14", "1", "0", "64.28")
    ,@(" Abu Dhabi", " October 07 2020", "KKR won by 10 runs", "Kolkata Knight Riders", "Chennai Super Kings", ("Nitish Rana" + $nbsp), "9", "10", "1", "0", "90.00")
)

$startRow = 16
$numericCols = @(7, 8, 9, 10, 11)   # G,H,I,J,K hold numeric-looking text; force text type

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($numericCols -contains $c) {
            # Without this, Excel auto-detects these digit/decimal strings
            # as numbers; the source data stores them as text.
            $cell.NumberFormat = "@"
            $cell.Value = $rowData[$c - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $rowData[$c - 1]
        }
    }
}
